# Update "Prix Spot" sheet (sheet 1): add a new day column AP ("25-jul")
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New header cell AP1 - copy formatting from the previous header cell (AO1)
# so the new cell keeps the exact same style (bold, centered, bordered).
$ws1.Range("AO1").Copy($ws1.Range("AP1")) | Out-Null
$ws1.Cells.Item(1, 42).Value2 = "25-jul"

$ap = @{
    2  = 88.69
    3  = 96.52
    4  = 78.48
    5  = 58.99
    6  = 51.97
    7  = 63.43
    8  = 77.97
    9  = 85.98
    10 = 98.34999999999999
    11 = 86.44
    12 = 62.4
    13 = 65.73999999999999
    14 = 59.01
    15 = 35.04
    16 = 16.2
    17 = 40.81
    18 = 46.49
    19 = 58.52
    20 = 63.23
    21 = 94.90000000000001
    22 = 85
    23 = 104.72
    24 = 105
    25 = 81.95999999999999
}

foreach ($row in $ap.Keys) {
    $ws1.Cells.Item($row, 42).Value2 = $ap[$row]
}

# Update "Gaz" sheet (sheet 2): append row 39 with new date/price
$ws2 = $wb.Worksheets.Item(2)
$a2 = $ws2.Cells.Item(39, 1)
# Build the new date as a formula first (so it is not auto-recognised as a
# real date value), then convert it to a plain static value - this keeps
# the cell a plain string like all the other date cells in the column.
$a2.Formula = "=""2025-07-23"""
$a2.Copy() | Out-Null
$a2.PasteSpecial(-4163) | Out-Null
$ws2.Cells.Item(39, 2).Value2 = 32.15

# Update "CO2" sheet (sheet 3): append row 39 with new date/price
$ws3 = $wb.Worksheets.Item(3)
$a3 = $ws3.Cells.Item(39, 1)
$a3.Formula = "=""2025-07-23"""
$a3.Copy() | Out-Null
$a3.PasteSpecial(-4163) | Out-Null
$ws3.Cells.Item(39, 2).Value2 = 68.40000000000001
